$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update totals for rows that changed before the deleted row
$ws.Range("C13").Value = 817
$ws.Range("C14").Value = 19
$ws.Range("C16").Value = 34
$ws.Range("C19").Value = 26
$ws.Range("C20").Value = 25

# Remove the "Maputo Provincia | Ponta De Ouro | 1" row entirely; rows below shift up
$ws.Rows.Item(21).Delete()

# After the deletion, former row 22 (Sofala | Beira | 168) is now row 21; update its total
$ws.Range("C21").Value = 164
